$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in End Date (J) values for existing rows 2, 3, 5 ---
$ws.Range("J2").Value = 43493
$ws.Range("J3").Value = 43493
$ws.Range("J5").Value = 43493

# Apply the same number format as the other date cells in column I/J
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# --- Add next two rows of "Actual Progression" data (rows 6 and 7) ---
$ws.Range("H6").Value = "Query coordination"
$ws.Range("I6").Value = 43493
$ws.Range("H7").Value = "Query optimization"
$ws.Range("I7").Value = 43493

# Apply matching date formatting to I6/I7 and the (still empty) J6/J7 cells
$ws.Range("I2").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J7").PasteSpecial(-4122)

# --- Expand Table2 (Actual Progression table) to include the new rows ---
$table2 = $ws.ListObjects.Item("Table2")
$table2.Resize($ws.Range("H1:J7"))

# --- Update the active cell selection to reflect where editing left off ---
$ws.Range("J11").Select()
